# Romania Liga I - daily odds-feed refresh (18-04-2024 00:36)
#
# Two things happened upstream:
#   1) The five fixtures that used to sit at rows 235-239 (ids 233-237,
#      all kicked off 2024-03-30 18:00) got re-sorted/re-fetched, so their
#      odds/result payloads now land on different rows than before - the
#      "id" column stays put but everything else (match id, teams, odds,
#      P/L columns...) cyclically rotates across those five rows.
#   2) Nine freshly scraped fixtures (ids 271-279) were appended at the
#      bottom of the sheet.
#
# This script reproduces both effects via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> 1-based column index, in sheet order.
$columns = @('A','B','C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC')
$colIndex = @{}
for ($i = 0; $i -lt $columns.Count; $i++) { $colIndex[$columns[$i]] = $i + 1 }

function Set-RowValues {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $columns) {
        if ($Values.ContainsKey($col)) {
            $ws.Cells.Item($Row, $colIndex[$col]).Value = $Values[$col]
        }
    }
}

function Copy-RowStyle {
    param(
        [int]$SrcRow,
        [int]$DstRow
    )
    # Only columns A (id, bold+border) and E (date, custom numfmt) carry an
    # explicit style in this sheet; every other column uses the default.
    $ws.Range("A$SrcRow").Copy() | Out-Null
    $ws.Range("A$DstRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("E$SrcRow").Copy() | Out-Null
    $ws.Range("E$DstRow").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# -----------------------------------------------------------------
# Part 1 - rows 235-239: same 5 fixtures, payload rotated by 2 rows
# (new row R shows what used to be on row R+2, wrapping within the block)
# -----------------------------------------------------------------
$rotated = @{
    235 = @{ "A"=233; "B"=6870268; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45359.625; "F"="Petrolul Ploiesti"; "G"="ACS Sepsi"; "H"=1; "I"=2; "J"="A"; "K"=2.8; "L"=3; "M"=2.55; "N"=3; "O"=3.2; "P"=2.3; "Q"=0.25; "R"=1.85; "S"=2; "T"=2.25; "U"=1.875; "V"=1.975; "W"=-1; "X"=-1; "Y"=1.3; "Z"=-1; "AA"=1; "AB"=0.875; "AC"=-1 };
    236 = @{ "A"=234; "B"=6865915; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45359.625; "F"="FC Voluntari"; "G"="Universitatea Cluj"; "H"=0; "I"=0; "J"="D"; "K"=3.5; "L"=3.25; "M"=2.05; "N"=3.4; "O"=3.1; "P"=2.15; "Q"=0.25; "R"=1.975; "S"=1.875; "T"=2.25; "U"=2.05; "V"=1.75; "W"=-1; "X"=2.1; "Y"=-1; "Z"=0.4875; "AA"=-0.5; "AB"=-1; "AC"=0.75 };
    237 = @{ "A"=235; "B"=6861095; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45359.625; "F"="FC Botosani"; "G"="Farul Constanta"; "H"=0; "I"=0; "J"="D"; "K"=3.75; "L"=3.4; "M"=1.909; "N"=3.1; "O"=3; "P"=2.375; "Q"=0.25; "R"=1.775; "S"=2.1; "T"=2; "U"=1.8; "V"=2.05; "W"=-1; "X"=2; "Y"=-1; "Z"=0.3875; "AA"=-0.5; "AB"=-1; "AC"=1.05 };
    238 = @{ "A"=236; "B"=6852370; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45359.625; "F"="Dinamo Bucharest"; "G"="ACS UTA Batrana Doamna"; "H"=1; "I"=0; "J"="H"; "K"=2.55; "L"=2.875; "M"=3; "N"=2.375; "O"=3; "P"=3.1; "Q"=-0.25; "R"=2; "S"=1.85; "T"=2.25; "U"=1.975; "V"=1.875; "W"=1.375; "X"=-1; "Y"=-1; "Z"=1; "AA"=-1; "AB"=-1; "AC"=0.875 };
    239 = @{ "A"=237; "B"=6836277; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45359.625; "F"="CFR Cluj"; "G"="AFC Hermannstadt"; "H"=1; "I"=0; "J"="H"; "K"=1.7; "L"=3.4; "M"=5; "N"=1.65; "O"=3.5; "P"=5.25; "Q"=-0.75; "R"=1.85; "S"=2; "T"=2.25; "U"=1.875; "V"=1.975; "W"=0.6499999999999999; "X"=-1; "Y"=-1; "Z"=0.425; "AA"=-0.5; "AB"=-1; "AC"=0.9750000000000001 }
}

foreach ($row in ($rotated.Keys | Sort-Object)) {
    Set-RowValues -Row $row -Values $rotated[$row]
}

# -----------------------------------------------------------------
# Part 2 - append 9 new fixtures as rows 273-281
# (rows 275-281 are still unplayed: no FTHG/FTAG/FTR/PL_AhOver/PL_AhUnder yet)
# -----------------------------------------------------------------
$newRows = @{
    273 = @{ "A"=271; "B"=7951754; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45396.625; "F"="CFR Cluj"; "G"="FCSB"; "H"=0; "I"=1; "J"="A"; "K"=2.3; "L"=3.1; "M"=3; "N"=2.15; "O"=3.1; "P"=3.4; "Q"=-0.25; "R"=1.875; "S"=1.975; "T"=2.25; "U"=1.875; "V"=1.975; "W"=-1; "X"=-1; "Y"=2.4; "Z"=-1; "AA"=0.9750000000000001; "AB"=-1; "AC"=0.9750000000000001 };
    274 = @{ "A"=272; "B"=7951756; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45397.60416666666; "F"="Rapid Bucuresti"; "G"="ACS Sepsi"; "H"=0; "I"=1; "J"="A"; "K"=1.95; "L"=3.25; "M"=3.75; "N"=1.666; "O"=3.6; "P"=4.75; "Q"=-0.75; "R"=1.875; "S"=1.975; "T"=2.75; "U"=1.95; "V"=1.9; "W"=-1; "X"=-1; "Y"=3.75; "Z"=-1; "AA"=0.9750000000000001; "AB"=-1; "AC"=0.8999999999999999 };
    275 = @{ "A"=273; "B"=7951793; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45401.47916666666; "F"="AFC Hermannstadt"; "G"="FC Botosani"; "K"=1.909; "L"=3; "M"=4.2; "N"=1.909; "O"=3; "P"=4.2; "Q"=-0.5; "R"=1.95; "S"=1.9; "T"=2; "U"=1.875; "V"=1.975; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    276 = @{ "A"=274; "B"=7951758; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45401.60416666666; "F"="Farul Constanta"; "G"="CFR Cluj"; "K"=3.1; "L"=3.2; "M"=2.2; "N"=3.1; "O"=3.2; "P"=2.2; "Q"=0.25; "R"=1.875; "S"=1.975; "T"=2.5; "U"=2.025; "V"=1.825; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    277 = @{ "A"=275; "B"=7951796; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45402.375; "F"="CSM Politehnica Iasi"; "G"="FC Voluntari"; "K"=2.25; "L"=3; "M"=3.25; "N"=2.25; "O"=3; "P"=3.25; "Q"=-0.25; "R"=1.975; "S"=1.875; "T"=2.25; "U"=2.05; "V"=1.8; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    278 = @{ "A"=276; "B"=7951797; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45402.47916666666; "F"="FC U Craiova 1948"; "G"="Dinamo Bucharest"; "K"=2.2; "L"=3.1; "M"=3.25; "N"=2.05; "O"=3.2; "P"=3.5; "Q"=-0.25; "R"=1.8; "S"=2.05; "T"=2.25; "U"=1.9; "V"=1.95; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    279 = @{ "A"=277; "B"=7951757; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45402.60416666666; "F"="FCSB"; "G"="Rapid Bucuresti"; "K"=1.85; "L"=3.5; "M"=4; "N"=1.85; "O"=3.5; "P"=4; "Q"=-0.5; "R"=1.9; "S"=1.95; "T"=2.5; "U"=2; "V"=1.85; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    280 = @{ "A"=278; "B"=7951795; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45403.41666666666; "F"="Otelul Galati"; "G"="Universitatea Cluj"; "K"=2.75; "L"=3.2; "M"=2.5; "N"=2.75; "O"=3.2; "P"=2.5; "Q"=0; "R"=2.025; "S"=1.825; "T"=2.25; "U"=1.975; "V"=1.875; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 };
    281 = @{ "A"=279; "B"=7951759; "C"="Romania Liga I"; "D"="Romania Liga I"; "E"=45403.54166666666; "F"="ACS Sepsi"; "G"="CS U Craiova"; "K"=2.6; "L"=3.2; "M"=2.6; "N"=2.55; "O"=3.2; "P"=2.625; "Q"=0; "R"=1.875; "S"=1.975; "T"=2.5; "U"=2.05; "V"=1.8; "W"=0; "X"=0; "Y"=0; "Z"=0; "AA"=0 }
}

$lastExistingRow = 272
foreach ($row in ($newRows.Keys | Sort-Object)) {
    Set-RowValues -Row $row -Values $newRows[$row]
    Copy-RowStyle -SrcRow $lastExistingRow -DstRow $row
}

$excel.CutCopyMode = $false
